$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1423
$ws.Range("I58").Value = 278.75
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 836.25
$ws.Range("L58").Value = 18000
$ws.Range("M58").Value = -686.25
$ws.Range("N58").Value = -18300

$ws.Range("H127").Value = 1616.409
$ws.Range("I127").Value = 1287.6
$ws.Range("J127").Value = 1890.4166
$ws.Range("K127").Value = 3862.8
$ws.Range("L127").Value = 5671.2498
$ws.Range("M127").Value = 1097.2
$ws.Range("N127").Value = -15591.2498

$ws.Range("H137").Value = 981.48
$ws.Range("I137").Value = 611.6667
$ws.Range("J137").Value = 1322.8462
$ws.Range("K137").Value = 1835.0001
$ws.Range("L137").Value = 3968.5386
$ws.Range("M137").Value = 714.9999
$ws.Range("N137").Value = -9068.5386

$ws.Range("H138").Value = 4025.8723
$ws.Range("I138").Value = 2270.0344
$ws.Range("J138").Value = 6854.722
$ws.Range("K138").Value = 6810.1032
$ws.Range("L138").Value = 20564.166
$ws.Range("M138").Value = -1670.1032
$ws.Range("N138").Value = -30844.166

$ws.Range("H141").Value = 3200.359
$ws.Range("I141").Value = 1422.6428
$ws.Range("K141").Value = 4267.928400000001
$ws.Range("M141").Value = 912.0715999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6119.2705
$ws.Range("I32").Value = 3780.9285
$ws.Range("J32").Value = 13394.111
$ws.Range("K32").Value = 3780.9285
$ws.Range("L32").Value = 13394.111
$ws.Range("M32").Value = -3493.9285
$ws.Range("N32").Value = -13968.111

$ws.Range("H61").Value = 1385.2285
$ws.Range("I61").Value = 701.15
$ws.Range("J61").Value = 2297.3333
$ws.Range("K61").Value = 701.15
$ws.Range("L61").Value = 2297.3333
$ws.Range("M61").Value = -489.15
$ws.Range("N61").Value = -2721.3333

$ws.Range("H74").Value = 987.44116
$ws.Range("I74").Value = 540.73914
$ws.Range("J74").Value = 1921.4546
$ws.Range("K74").Value = 540.73914
$ws.Range("L74").Value = 1921.4546
$ws.Range("M74").Value = 333.26086
$ws.Range("N74").Value = -3669.4546

$ws.Range("H77").Value = 987.44116
$ws.Range("I77").Value = 540.73914
$ws.Range("J77").Value = 1921.4546
$ws.Range("K77").Value = 2703.6957
$ws.Range("L77").Value = 9607.273000000001
$ws.Range("M77").Value = 1664.3043
$ws.Range("N77").Value = -18343.273

$ws.Range("H136").Value = 1385.2285
$ws.Range("I136").Value = 701.15
$ws.Range("J136").Value = 2297.3333
$ws.Range("K136").Value = 2103.45
$ws.Range("L136").Value = 6891.999899999999
$ws.Range("M136").Value = 446.5500000000002
$ws.Range("N136").Value = -11991.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7937998.5
$ws.Range("I31").Value = 21740062
$ws.Range("J31").Value = 1812.1
$ws.Range("K31").Value = 21740062
$ws.Range("L31").Value = 1812.1
$ws.Range("M31").Value = -21739767
$ws.Range("N31").Value = -2402.1

$ws.Range("H34").Value = 7937998.5
$ws.Range("I34").Value = 21740062
$ws.Range("J34").Value = 1812.1
$ws.Range("K34").Value = 21740062
$ws.Range("L34").Value = 1812.1
$ws.Range("M34").Value = -21739860
$ws.Range("N34").Value = -2216.1

$ws.Range("H58").Value = 2110.0588
$ws.Range("J58").Value = 2137.25
$ws.Range("L58").Value = 2137.25
$ws.Range("N58").Value = -2543.25

$ws.Range("H99").Value = 125000670
$ws.Range("I99").Value = 793
$ws.Range("J99").Value = 500000300
$ws.Range("K99").Value = 793
$ws.Range("L99").Value = 500000300
$ws.Range("M99").Value = 705
$ws.Range("N99").Value = -500003296

$ws.Range("H122").Value = 1340
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -9400

$ws.Range("H126").Value = 125000670
$ws.Range("I126").Value = 793
$ws.Range("J126").Value = 500000300
$ws.Range("K126").Value = 2379
$ws.Range("L126").Value = 1500000900
$ws.Range("M126").Value = 91
$ws.Range("N126").Value = -1500005840

$ws.Range("H132").Value = 3806.5264
$ws.Range("I132").Value = 3643.8333
$ws.Range("J132").Value = 4085.4285
$ws.Range("K132").Value = 10931.4999
$ws.Range("L132").Value = 12256.2855
$ws.Range("M132").Value = -8401.499899999999
$ws.Range("N132").Value = -17316.2855

$ws.Range("H134").Value = 15626266
$ws.Range("I134").Value = 20001100
$ws.Range("K134").Value = 60003300
$ws.Range("M134").Value = -60000765

$ws.Range("H136").Value = 2110.0588
$ws.Range("J136").Value = 2137.25
$ws.Range("L136").Value = 6411.75
$ws.Range("N136").Value = -11511.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 337.45834
$ws.Range("I107").Value = 242.07812
$ws.Range("J107").Value = 1100.5
$ws.Range("K107").Value = 726.23436
$ws.Range("L107").Value = 3301.5
$ws.Range("M107").Value = 1193.76564
$ws.Range("N107").Value = -7141.5

$ws.Range("H118").Value = 1739298.2
$ws.Range("I118").Value = 2249.25
$ws.Range("J118").Value = 2318314.8
$ws.Range("K118").Value = 6747.75
$ws.Range("L118").Value = 6954944.399999999
$ws.Range("M118").Value = -5504.75
$ws.Range("N118").Value = -6957430.399999999

$ws.Range("H131").Value = 10639207
$ws.Range("I131").Value = 41667150
$ws.Range("J131").Value = 1053.1714
$ws.Range("K131").Value = 125001450
$ws.Range("L131").Value = 3159.5142
$ws.Range("M131").Value = -124996410
$ws.Range("N131").Value = -13239.5142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3220.25
$ws.Range("I122").Value = 2603.182
$ws.Range("K122").Value = 7809.545999999999
$ws.Range("M122").Value = -5359.545999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2579.75
$ws.Range("I40").Value = 2356.4443
$ws.Range("J40").Value = 3249.6667
$ws.Range("K40").Value = 2356.4443
$ws.Range("L40").Value = 3249.6667
$ws.Range("M40").Value = -2220.4443
$ws.Range("N40").Value = -3521.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16251.429
$ws.Range("I81").Value = 17293.5
$ws.Range("J81").Value = 9999
$ws.Range("K81").Value = 34587
$ws.Range("L81").Value = 19998
$ws.Range("M81").Value = -33526
$ws.Range("N81").Value = -22120

$ws.Range("H84").Value = 16251.429
$ws.Range("I84").Value = 17293.5
$ws.Range("J84").Value = 9999
$ws.Range("K84").Value = 172935
$ws.Range("L84").Value = 99990
$ws.Range("M84").Value = -167631
$ws.Range("N84").Value = -110598

$ws.Range("H136").Value = 13754.333
$ws.Range("I136").Value = 3928
$ws.Range("J136").Value = 18667.5
$ws.Range("K136").Value = 11784
$ws.Range("L136").Value = 56002.5
$ws.Range("M136").Value = -9234
$ws.Range("N136").Value = -61102.5
